$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from J1 onto K1 so the new header cell matches
# the existing bold/bordered/centered header formatting (style index 1).
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# intervention_type values for rows 1 (header) through 114.
# Empty strings correspond to rows where the sponsor left the field blank.
$values = @(
    "intervention_type",
    "DRUG",
    "DRUG",
    "OTHER",
    "OTHER",
    "DRUG",
    "DRUG",
    "OTHER",
    "DRUG",
    "PROCEDURE",
    "DEVICE",
    "DRUG",
    "DRUG",
    "DRUG",
    "PROCEDURE",
    "DEVICE",
    "PROCEDURE",
    "BEHAVIORAL",
    "DRUG",
    "BEHAVIORAL",
    "PROCEDURE",
    "DRUG",
    "BIOLOGICAL",
    "DEVICE",
    "PROCEDURE",
    "OTHER",
    "DRUG",
    "OTHER",
    "OTHER",
    "OTHER",
    "BIOLOGICAL",
    "DRUG",
    "PROCEDURE",
    "DEVICE",
    "PROCEDURE",
    "BIOLOGICAL",
    "BEHAVIORAL",
    "OTHER",
    "OTHER",
    "DRUG",
    "DRUG",
    "PROCEDURE",
    "DEVICE",
    "DEVICE",
    "DEVICE",
    "DEVICE",
    "DRUG",
    "DEVICE",
    "",
    "DRUG",
    "DEVICE",
    "DEVICE",
    "DIAGNOSTIC_TEST",
    "OTHER",
    "OTHER",
    "OTHER",
    "OTHER",
    "DRUG",
    "DEVICE",
    "BIOLOGICAL",
    "DRUG",
    "OTHER",
    "BEHAVIORAL",
    "DIAGNOSTIC_TEST",
    "OTHER",
    "DRUG",
    "DEVICE",
    "OTHER",
    "OTHER",
    "OTHER",
    "OTHER",
    "OTHER",
    "DRUG",
    "BIOLOGICAL",
    "OTHER",
    "OTHER",
    "OTHER",
    "OTHER",
    "BIOLOGICAL",
    "OTHER",
    "PROCEDURE",
    "PROCEDURE",
    "PROCEDURE",
    "BIOLOGICAL",
    "OTHER",
    "DRUG",
    "OTHER",
    "OTHER",
    "DIAGNOSTIC_TEST",
    "PROCEDURE",
    "OTHER",
    "OTHER",
    "DRUG",
    "PROCEDURE",
    "BIOLOGICAL",
    "BIOLOGICAL",
    "OTHER",
    "OTHER",
    "OTHER",
    "OTHER",
    "DIAGNOSTIC_TEST",
    "OTHER",
    "OTHER",
    "BIOLOGICAL",
    "BEHAVIORAL",
    "DRUG",
    "PROCEDURE",
    "OTHER",
    "OTHER",
    "PROCEDURE",
    "DRUG",
    "",
    "",
    ""
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $val = $values[$i]
    if ($val -ne "") {
        $ws.Cells.Item($row, 11).Value = $val
    }
}